$d = $word.ActiveDocument

# Locate the paragraph that contains the unique marker text for the
# "https://www.youtube.com/watch?v=70LxGwhWK8U" link - this is the
# paragraph whose formatting (pPr mark + run) needs the accent3 color.
$rng = $d.Content
$found = $rng.Find.Execute("70LxGwhWK8U", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target paragraph text"
}

$para = $rng.Paragraphs(1).Range

# Pull this paragraph's canonical OOXML so we can surgically add the
# missing <w:color/> element to both the paragraph-mark run properties
# and the text run properties, while preserving every other attribute
# (rsid values, etc.) unchanged.
$xml = $para.WordOpenXML

$startTag = "<w:p "
$startIdx = $xml.IndexOf($startTag)
$endTag = "</w:p>"
$endIdx = $xml.IndexOf($endTag) + $endTag.Length
$paraXml = $xml.Substring($startIdx, $endIdx - $startIdx)

# Strip the synthetic w14:paraId / w14:textId attributes that the
# round-trip serializer injects; the original document does not use
# those attributes on this paragraph.
$paraXml = $paraXml -replace ' w14:paraId="[0-9A-Fa-f]+"', ''
$paraXml = $paraXml -replace ' w14:textId="[0-9A-Fa-f]+"', ''

# Insert the accent3 theme color right after the rFonts element inside
# every run-properties block that carries the majorHAnsi font triplet
# used by this paragraph (covers both the paragraph mark's rPr and the
# run's rPr).
$colorElem = '<w:color w:val="9BBB59" w:themeColor="accent3"/>'
$needle = '<w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>'
if ($paraXml.IndexOf($colorElem) -lt 0) {
    $paraXml = $paraXml.Replace($needle, $needle + $colorElem)
}

$pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $paraXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$null = $para.InsertXML($pkg)
